$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A5").NumberFormat = "@"

$ws.Range("A2").Value = "110007405"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 753

$ws.Range("A3").Value = "111937242"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 300

$ws.Range("A4").Value = "111904125"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 300

$ws.Range("A5").Value = "7114168"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 753

$ws.Range("C6").Value = "Total: 2106"
